$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the leading "+1" country code from the existing phone numbers.
# (Leading apostrophe keeps the COM layer from re-inferring these as numbers,
# matching the workbook's existing text-formatted phone column.)
$ws.Range("B2").Value = "'8733083608"
$ws.Range("B3").Value = "'9753501446"
$ws.Range("B4").Value = "'3544172164"
$ws.Range("B5").Value = "'7713556101"

# Append the newly-merged rows of contact data.
$newRows = @(
    @("jane.taylor@example.com",   "8733083608", "Jane Taylor"),
    @("daniel.williams@example.com","9753501446", "Daniel Williams"),
    @("jane.williams@example.com", "3544172164", "Jane Williams"),
    @("daniel.williams@example.com","6646753997", "Daniel Williams"),
    @("alice.williams@example.com","7713556101", "Alice Williams"),
    @("jane.davis@example.com",    "9123981617", "Jane Davis"),
    @("jane.davis@example.com",    "9123981618", "Jane Davis")
)

$r = 6
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Trailing separator row, shaded light grey.
$ws.Range("A13:E13").Interior.Color = 0xD9D9D9
